$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.672.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.216.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.13"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.66%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.73%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.546.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.897"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.223.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.671.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +22.24%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.00%  "
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0745"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.34%  "
$ws.Range("E38").Value = "  +4.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0303"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +27.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.203"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.46%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  +5.81%  "
